$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text cells keep their exact string representation (avoid Excel
# auto-converting numeric-looking strings like "14.20" or "1.00" into numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.608.41'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.16%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.965.91'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '243.68'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.621'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.14%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.26'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.14%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +2.67%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -3.03%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.20'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.845'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.19%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.255.49'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.66'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.29'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.963.41'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.50%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.527.23'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.40'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '229.43'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.76%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.45'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.144'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +7.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.15'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.38%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.99'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.31'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.32'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +19.06%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.48%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +2.60%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.22%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.55'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.84%  '
$ws.Range('B35').NumberFormat = '@'
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').NumberFormat = '@'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('B36').NumberFormat = '@'
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').NumberFormat = '@'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +2.73%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.56%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -13.19%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0970'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.94%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.93'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -0.75%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.371.74'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.84'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.06%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +1.04%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.43%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.11'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.89%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.150.33'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.95%  '
